# Remove the "Transition_Name_Annot" and "ISTD_Annot" worksheets, leaving
# only "Sample_Annot" in the workbook (matching the new MSTemplate_Creator
# output format for the Sample Annotation validation unit test).

$wb = $excel.ActiveWorkbook

# Turn off alerts so sheet deletion doesn't prompt for confirmation.
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Transition_Name_Annot").Delete()
$wb.Worksheets.Item("ISTD_Annot").Delete()

$excel.DisplayAlerts = $true

# Make sure the remaining sheet is active/selected.
$wb.Worksheets.Item("Sample_Annot").Activate()
